$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Add the new "WMT_Extract_SA" worksheet as the last tab (after "T2A"),
#    matching the workbook.xml <sheets> / activeTab changes in the diff.
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$newSheet.Name = "WMT_Extract_SA"

# ---------------------------------------------------------------------------
# 2. Populate the header row with the column names used by the new sheet.
#    These reuse existing shared strings where already present in the
#    workbook, and add three brand-new ones (Disposal_Type_Desc,
#    Disposal_Type_Code, Standalone_Order).
# ---------------------------------------------------------------------------
$newSheet.Range("A1").Value = "Case_Ref_No"
$newSheet.Range("B1").Value = "Tier_Code"
$newSheet.Range("C1").Value = "Team_Code"
$newSheet.Range("D1").Value = "OM_Grade_Code"
$newSheet.Range("E1").Value = "OM_Key"
$newSheet.Range("F1").Value = "Location"
$newSheet.Range("G1").Value = "Disposal_Type_Desc"
$newSheet.Range("H1").Value = "Disposal_Type_Code"
$newSheet.Range("I1").Value = "Standalone_Order"

# ---------------------------------------------------------------------------
# 3. Apply the header formatting. Start from the existing dark-header style
#    used elsewhere in the workbook (numFmt "@", Arial 9 bold white on
#    #666699) by pasting formats from Court_Reports!C1, then recolor the
#    borders to the new light lavender (#CCCCFF) used for this sheet -
#    full box on A1, no left edge on the rest (to avoid doubled borders).
# ---------------------------------------------------------------------------
$courtReports = $wb.Worksheets.Item("Court_Reports")
$courtReports.Range("C1").Copy()
$newSheet.Range("A1:I1").PasteSpecial(-4122)

$cellA1 = $newSheet.Range("A1")
$cellRest = $newSheet.Range("B1:I1")

$cellA1.Borders.Color = 16764108

$cellRest.Borders.Color = 16764108
foreach ($cell in $cellRest) {
    $cell.Borders.Item(7).LineStyle = -4142
}

# Select the header row on the new sheet (it becomes the active/visible tab).
$newSheet.Range("A1:I1").Select()

Write-Output "Added WMT_Extract_SA worksheet with header row"
